$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix 1: "distitntas" -> "distintas"
# The misspelled word sits inside its own w:proofErr spellStart/spellEnd pair.
# Replacing text that crosses BOTH proofErr markers collapses the block into a
# single run and drops the (now unneeded) proofErr wrapper, matching the
# post-edit XML where "distintas" is no longer flagged.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("tienes distitntas fases iniciando", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$r.Text = "tienes distintas fases iniciando"

# Re-split the merged run back into the three pieces the target keeps separate:
#   "...tienes "  |  "distintas"  |  " fases iniciando en \u201c"
$rWord = $d.Content
$rWord.Find.Execute("distintas", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$rWord.Bold = 1
$rWord.Bold = 0

$rTail = $d.Content
$rTail.Find.Execute(" fases iniciando en $([char]0x201C)", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$rTail.Bold = 1
$rTail.Bold = 0

# ---------------------------------------------------------------------------
# Fix 2: merge the separate " " and "introducido" runs into one run
# " introducido", while keeping the following runs (" el dinero", ", ",
# "pasamos al estado final \u201c") split exactly as before.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("introducido", $true, $false, $false, $false, $false,
                  $true, 1, $false, "introducido", 2)

# Re-split the merged tail back into its original pieces, working right-to-left
# so each boundary search is unambiguous:
#   " introducido" | " el dinero" | ", " | "pasamos al estado final \u201c"
$s1 = $d.Content
$s1.Find.Execute("pasamos al estado final $([char]0x201C)", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0)
$s1.Bold = 1
$s1.Bold = 0

$s2 = $d.Content
$s2.Find.Execute(", ", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0)
$s2.Bold = 1
$s2.Bold = 0

$s3 = $d.Content
$s3.Find.Execute(" el dinero", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0)
$s3.Bold = 1
$s3.Bold = 0
